$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New "Data" column (E) values for years 1950-2016 (row -> text value),
# keyed by worksheet row number. Rows 2-60 (years 1950-2008) get updated
# values; rows 61-68 (years 2009-2016) are new rows appended below.
$dataValues = @{
    2 = "982"
    3 = "1015"
    4 = "1022"
    5 = "1036"
    6 = "1068"
    7 = "1074"
    8 = "1086"
    9 = "1097"
    10 = "1097"
    11 = "1129"
    12 = "1167"
    13 = "1318"
    14 = "1262"
    15 = "1232"
    16 = "1274"
    17 = "1368"
    18 = "1545"
    19 = "1548"
    20 = "1519"
    21 = "1629"
    22 = "1422"
    23 = "1376"
    24 = "1334"
    25 = "1411"
    26 = "1452"
    27 = "1497"
    28 = "1436"
    29 = "1417"
    30 = "1446"
    31 = "1360"
    32 = "1325"
    33 = "1521"
    34 = "1581"
    35 = "1482"
    36 = "1462"
    37 = "1451"
    38 = "1422"
    39 = "1406"
    40 = "1417"
    41 = "1407"
    42 = "1400"
    43 = "1390.70497512334"
    44 = "1391.40545449168"
    45 = "1385.5908111508"
    46 = "1345.0349874729"
    47 = "1317.24249573015"
    48 = "1336.84732412669"
    49 = "1308.48667194314"
    50 = "1324.92449635651"
    51 = "1385.63177687098"
    52 = "1440.86192167419"
    53 = "1495.091836157"
    54 = "1473.40067435036"
    55 = "1478.95622846008"
    56 = "1585.79921991557"
    57 = "1541.61159787083"
    58 = "1530.68440032576"
    59 = "1557.84736603492"
    60 = "1618.44448244572"
    61 = "1692.64234911445"
    62 = "1772.72685288389"
    63 = "1669"
    64 = "1725"
    65 = "1767"
    66 = "1743"
    67 = "1786"
    68 = "1788"
}

# Append the new rows (61-68) for years 2009-2016: Country Code, Country
# Name, Indicator and Year columns follow the same pattern as existing rows.
for ($row = 61; $row -le 68; $row++) {
    $year = 1950 + ($row - 2)
    $ws.Cells.Item($row, 1).Value = 270
    $ws.Cells.Item($row, 2).Value = "The Gambia"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $year
}

# Write column E (Data) as text for every row, rows 2-68, in one pass.
# NumberFormat is forced to text first (so numeric-looking strings like
# "982" aren't auto-converted to numbers by Excel), then cleared again
# afterwards so no extra cell formatting is left behind.
$dataRange = $ws.Range("E2:E68")
$dataRange.NumberFormat = "@"
foreach ($row in $dataValues.Keys) {
    $ws.Cells.Item($row, 5).Value = $dataValues[$row]
}
$dataRange.ClearFormats()
